$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")
$ws2.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
Write-Output "scrolled"
